$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 259. This shifts the existing rows
# 259-357 down to 260-358 (dimension grows from A1:R357 to A1:R358), then
# we populate the newly inserted row with the new record's values.
$ws.Rows.Item(259).Insert()

$ws.Cells.Item(259, 1).Value = 5
$ws.Cells.Item(259, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(259, 3).Value = "Maule"
$ws.Cells.Item(259, 4).Value = 44755
$ws.Cells.Item(259, 5).Value = 7
$ws.Cells.Item(259, 6).Value = 100112023
$ws.Cells.Item(259, 7).Value = "Brócoli"
$ws.Cells.Item(259, 8).Value = "Sin especificar"
$ws.Cells.Item(259, 9).Value = "Primera"
$ws.Cells.Item(259, 10).Value = 4000
$ws.Cells.Item(259, 11).Value = 800
$ws.Cells.Item(259, 12).Value = 800
$ws.Cells.Item(259, 13).Value = 800
$ws.Cells.Item(259, 14).Value = '$/unidad'
$ws.Cells.Item(259, 15).Value = "Región del Maule"
$ws.Cells.Item(259, 16).Value = 800
$ws.Cells.Item(259, 17).Value = 1
$ws.Cells.Item(259, 18).Value = "Hortaliza"
